$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New search_radius (boundary) value applied to rows 2-13
$searchRadius = 135000

# VAR(OK) values (column D)
$varOK = @(
    0.2546468043688802,
    0.4110409802702994,
    0.5355094759621808,
    0.56597042598376,
    0.5794681151118293,
    0.5872462595146157,
    0.6063895716669161,
    0.6185023987743629,
    0.638018318454436,
    0.6578425679565504,
    0.6648155393024699,
    0.6673271402185389
)

# MSPE values (column E)
$mspe = @(
    4.090477886565919,
    3.912426298598175,
    3.754207723738251,
    3.725360916877521,
    3.706315775055423,
    3.69561847722617,
    3.672848038905631,
    3.662001846739384,
    3.642633110237823,
    3.650818488778145,
    3.651828017178621,
    3.653356167919833
)

# S_nugget values (column F)
$sNugget = @(
    1.045,
    1.6292,
    2.0564,
    2.4922,
    2.7495,
    2.9962,
    3.1795,
    3.3527,
    3.4535,
    3.4535,
    3.4535,
    3.4535
)

# VAR(DATA) values (column H)
$varData = @(
    4.246647459488038,
    4.246647459488038,
    4.246647459488038,
    4.246647459488038,
    4.246647459488038,
    4.246647459488038,
    4.246647459488038,
    4.246647459488038,
    4.246647459488038,
    4.246647459488038,
    4.246647459488038,
    4.246647459488038
)

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $searchRadius
    $ws.Cells.Item($row, 4).Value = $varOK[$i]
    $ws.Cells.Item($row, 5).Value = $mspe[$i]
    $ws.Cells.Item($row, 6).Value = $sNugget[$i]
    $ws.Cells.Item($row, 8).Value = $varData[$i]
}
